$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in "Razon social" / "Nombre Fantasia" entries where a
# --- trailing/inter-name comma was scraped as a literal comma instead of a
# --- period (and a trailing "S.H." normalized to "SH").
$ws.Range("E85").Value = "URUMAT SOCIEDAD SIMPLE DE BONASEGLA CATALINA. BONASEGLA LUCIANA Y BONASEGLA SILVIO"
$ws.Range("E107").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E219").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E232").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("F294").Value = "CLERICE. MIGUEL ANGEL"

# --- Fix formatting of the "Importe" column: values were scraped using the
# --- Argentine locale convention ("." thousands separator, "," decimal
# --- separator, e.g. "126.600,00") and must be normalized to a plain
# --- decimal representation ("126600.00").
$importe = @(
    "126600.00", "50400.00", "15600.00", "20800.00", "76232.50", "298.00", "2375.00", "2100.00", "2500.00", "387000.00",
    "263828.90", "82500.00", "809000.00", "690.00", "9800.00", "420.00", "958296.75", "1312993.97", "103712.00", "102300.00",
    "98983.80", "222941.72", "192081.69", "5600.00", "16000.00", "3084.00", "20400.00", "19291.60", "148562.96", "11239.00",
    "84007.16", "53275.00", "9750.00", "1500.00", "13500.00", "20960.00", "42000.32", "11400.00", "6880.00", "1367.98",
    "296891.10", "1020.22", "56.00", "407.00", "403148.47", "333254.90", "500.00", "2400.00", "25917.66", "236100.00",
    "80.00", "684400.00", "1080.00", "341204.96", "340.00", "1998716.51", "13325.00", "7452.14", "1366.25", "22343.00",
    "97881.53", "280.00", "400.00", "16205.86", "312.20", "575.00", "2300.00", "12971.74", "2376.60", "958.00",
    "3401.08", "9659.00", "1980.00", "10300.00", "3435.00", "14400.00", "14182.40", "3529.75", "9611.10", "3090.00",
    "1442.72", "3578.67", "450.00", "13500.00", "110927.00", "19515.53", "48100.00", "35600.00", "60400.00", "42000.00",
    "73500.00", "87500.00", "4527.15", "352.00", "5358.46", "300.00", "12535.00", "650.00", "40260.00", "690.00",
    "151500.00", "22300.00", "4480.00", "930.00", "45216.00", "1800.00", "4520.00", "6321.00", "200.00", "23.04",
    "31652.08", "284.83", "39587.19", "23500.46", "215.06", "849.00", "12535.00", "5026.02", "17405.00", "13999.00",
    "1348.50", "190.00", "550.00", "2769.00", "8220.00", "12129.00", "10395.00", "4600.00", "73584.00", "15673.00",
    "95211.92", "2199.00", "2680.00", "13569.96", "400.00", "1491.70", "68685.00", "1677.50", "4160.00", "23552.00",
    "14080.00", "69641.19", "1810.40", "325184.00", "58700.00", "8600.00", "9000.00", "15800.00", "9368.00", "126280.00",
    "10848.00", "11000.00", "3000.00", "92500.00", "8500.00", "17000.00", "99900.00", "30000.00", "19000.00", "35000.00",
    "161000.00", "12000.00", "40000.00", "29000.00", "18811.49", "20126.00", "4387.21", "144890.00", "8245.00", "552.64",
    "1380.50", "1191.28", "12485.00", "35805.00", "9000.00", "16500.00", "12000.00", "22000.00", "7000.00", "24000.00",
    "10000.00", "8500.00", "12000.00", "20000.00", "10000.00", "9000.00", "9000.00", "4000.00", "20000.00", "10000.00",
    "5000.00", "21000.00", "10000.00", "10000.00", "10000.00", "5000.00", "16000.00", "10000.00", "12000.00", "4500.00",
    "5000.00", "13500.00", "10500.00", "48000.00", "18000.00", "10000.00", "10000.00", "113875.55", "7000.00", "35000.00",
    "10000.00", "4810.00", "37300.00", "57500.00", "25800.00", "55686.30", "12400.00", "1500.00", "8900.00", "2527.24",
    "3030.02", "50390.00", "21100.00", "879.96", "6880.00", "1510.00", "8600.00", "7500.00", "13560.00", "31160.17",
    "23270.00", "4890.00", "3006.02", "2871.00", "2147.28", "857.80", "372680.00", "1019.96", "569.77", "310.00",
    "33936.27", "10817.93", "80000.00", "40000.00", "40000.00", "40000.00", "80000.00", "40000.00", "55000.00", "40000.00",
    "40000.00", "80000.00", "80000.00", "79500.00", "9000.00", "42831.00", "6024925.03", "6900.00", "111800.00", "61180.00",
    "114800.00", "70889.00", "3965000.00", "241000.00", "296700.00", "241000.00", "241000.00", "241000.00", "245950.00", "451000.00",
    "241000.00", "573450.00", "513000.00", "306700.00", "241000.00", "241000.00", "482000.00", "419200.00", "452200.00", "697450.00",
    "451000.00", "739100.00", "482000.00", "247350.00", "205702.00", "763049.41", "1534876.54", "459000.00", "106590.00", "2163489.78",
    "1399250.00", "15000.00", "8000.00", "118000.00", "7500.00", "60000.00", "10500.00", "6500.00", "176500.00", "6500.00",
    "20000.00", "5917.94", "8000.00", "2600.00", "1090.00", "318170.00", "335000.00"
)

$startRow = 2
for ($i = 0; $i -lt $importe.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $importe[$i]
}
